# Adding feature selection script
# Add a new "name" column (N) with header and a single data value in row 3,
# mirroring the formatting of the existing header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell N1: bold + centered, same as the rest of row 1
$ws.Range("N1").Value = "name"
$ws.Range("N1").Font.Bold = $true
$ws.Range("N1").HorizontalAlignment = -4108  # xlCenter

# New data cell N3: plain formatting, same as the rest of row 3
$ws.Range("N3").Value = "Overall QOL transformed Score"
